$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; this shifts the existing D:K data block one column
# right to E:L (matching the new <dimension ref="A5:L102"/>) and carries each cells
# original value/type along with it.
$ws.Columns("D:D").Insert(-4161, 0)

# The newly inserted column D has no number format yet. Copy the formatting from column E
# (the old column D, now one column over) into D so the new column keeps the same date /
# number formatting as the rest of the table. Done per contiguous data block (Income
# Statement, Balance Sheet, Cash Flow Statement) so the blank separator rows (36/37, 78/79)
# are left untouched.
$ws.Range("E7:E35").Copy() | Out-Null
$ws.Range("D7:D35").PasteSpecial(-4122) | Out-Null
$ws.Range("E38:E77").Copy() | Out-Null
$ws.Range("D38:D77").PasteSpecial(-4122) | Out-Null
$ws.Range("E80:E102").Copy() | Out-Null
$ws.Range("D80:D102").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Populate the new periods figures in column D, plus a few restated prior-period values
# in column F (the old column E) that shifted right along with everything else.
$ws.Range("D7").Value = 43465
$ws.Range("D8").Value = 4131200
$ws.Range("F8").Value = 3084100
$ws.Range("D9").Value = 374300
$ws.Range("F9").Value = 258300
$ws.Range("D10").Value = 3756900
$ws.Range("F10").Value = 2825800
$ws.Range("D12").Value = 730400
$ws.Range("F12").Value = 757200
$ws.Range("D13").Value = 0
$ws.Range("D14").Value = 1208500
$ws.Range("F14").Value = 90300
$ws.Range("D15").Value = 320100
$ws.Range("F15").Value = 322200
$ws.Range("D17").Value = 3861600
$ws.Range("F17").Value = 2416700
$ws.Range("D18").Value = 269600
$ws.Range("F18").Value = 667400
$ws.Range("D20").Value = 70800
$ws.Range("F20").Value = 5700
$ws.Range("D21").Value = 745700
$ws.Range("F21").Value = 1069500
$ws.Range("D22").Value = 98200
$ws.Range("F22").Value = 96900
$ws.Range("D23").Value = 242200
$ws.Range("F23").Value = 576200
$ws.Range("D24").Value = 221100
$ws.Range("F24").Value = 176800
$ws.Range("D25").Value = 0
$ws.Range("D26").Value = 21100
$ws.Range("F26").Value = 399400
$ws.Range("D27").Value = 21100
$ws.Range("F27").Value = 399400
$ws.Range("D28").Value = 0
$ws.Range("D29").Value = 56500
$ws.Range("D30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("D32").Value = -70800
$ws.Range("F32").Value = -5700
$ws.Range("D33").Value = 77600
$ws.Range("F33").Value = 399400
$ws.Range("D34").Value = 0
$ws.Range("D35").Value = 77600
$ws.Range("F35").Value = 399400
$ws.Range("D38").Value = 43465
$ws.Range("D41").Value = 1365500
$ws.Range("D42").Value = 198300
$ws.Range("D43").Value = 922300
$ws.Range("D44").Value = 472500
$ws.Range("D45").Value = 426400
$ws.Range("D46").Value = 3385000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("F47").Value = 0
$ws.Range("G47").Value = 0
$ws.Range("H47").Value = 0
$ws.Range("I47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("D48").Value = 1471500
$ws.Range("D49").Value = 8678700
$ws.Range("D50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("D52").Value = 396700
$ws.Range("D53").Value = 0
$ws.Range("D54").Value = 13931900
$ws.Range("D57").Value = 74400
$ws.Range("D58").Value = 343800
$ws.Range("D59").Value = 755800
$ws.Range("D60").Value = 1174000
$ws.Range("D61").Value = 2862700
$ws.Range("D62").Value = 729900
$ws.Range("D63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("D66").Value = 4766600
$ws.Range("D68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("D72").Value = 2325800
$ws.Range("D73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("D76").Value = 9165300
$ws.Range("D77").Value = 0
$ws.Range("D80").Value = 43465
$ws.Range("D81").Value = 77600
$ws.Range("F81").Value = 399400
$ws.Range("D83").Value = 405300
$ws.Range("F83").Value = 396400
$ws.Range("D84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("D89").Value = 426000
$ws.Range("F89").Value = 1086300
$ws.Range("D91").Value = -213000
$ws.Range("F91").Value = -332700
$ws.Range("D92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("D94").Value = 470500
$ws.Range("F94").Value = -287600
$ws.Range("D96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("D100").Value = -102400
$ws.Range("F100").Value = -836200
$ws.Range("D101").Value = -11200
$ws.Range("F101").Value = -6600
$ws.Range("D102").Value = 782900
$ws.Range("F102").Value = -44100
